$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: column C "Value" -> "Participation"
$ws.Cells.Item(1, 3).Value = "Participation"

# New participation counts for years 1980-2023 (rows 2-45), replacing the
# old dollar-value figures in column C.
$participation = @{
    1980 = 54;  1981 = 61;  1982 = 51;  1983 = 44;  1984 = 7;
    1985 = 45;  1986 = 43;  1987 = 41;  1988 = 51;  1989 = 48;
    1990 = 43;  1991 = 37;  1992 = 18;  1993 = 44;  1994 = 42;
    1995 = 59;  1996 = 64;  1997 = 56;  1998 = 20;  1999 = 76;
    2000 = 88;  2001 = 62;  2002 = 73;  2003 = 43;  2004 = 73;
    2005 = 90;  2006 = 89;  2007 = 61;  2008 = 80;  2009 = 76;
    2010 = 76;  2011 = 79;  2012 = 97;  2013 = 61;  2014 = 36;
    2015 = 12;  2016 = 51;  2017 = 49;  2018 = 36;  2019 = 29;
    2020 = 31;  2021 = 39;  2022 = 36;  2023 = 45
}

for ($row = 2; $row -le 45; $row++) {
    $year = $ws.Cells.Item($row, 1).Value()
    $newVal = $participation[[int]$year]
    $ws.Cells.Item($row, 3).Value = $newVal
}

# Total row (46): column C total becomes the sum of the new participation
# counts; restyle it to match column B's total style (plain number, bold).
$ws.Cells.Item(46, 3).Value = 504

# Column C (rows 2-46) should now share the same number format as column B
# (plain thousands separator, no currency sign) instead of the old "$" format.
$ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item(46, 3)).NumberFormat = "#,0"
